$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column S (rows 4-14) into column T so the new
# column inherits identical cell styles (reuses existing style records).
for ($r = 4; $r -le 14; $r++) {
    $ws.Cells.Item($r, 19).Copy()
    $ws.Cells.Item($r, 20).PasteSpecial(-4122)  # xlPasteFormats
}

# Add new data values for column T (year 2023)
$ws.Range("T4").Value = 2023
$values = @{
    "T5"  = 4.8187602774004432
    "T6"  = 11.788953009068425
    "T7"  = 5.2855407047387608
    "T8"  = 11.35112240576027
    "T9"  = 16.577540106951872
    "T10" = 14.651002073255009
    "T11" = 5.034965034965035
    "T12" = 3.1837160751565765
    "T13" = 2.2263731825525039
    "T14" = 5.1321450522433931
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Reset the selection away from the stale "T6" reference left in the
# source file (the published workbook is re-generated without a
# meaningful user selection, so park the cursor back at A1).
$ws.Range("A1").Select()
